$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("I2").Value = 0
$ws1.Range("M2").Value = 0
$ws1.Range("Q2").Value = 0
$ws1.Range("L3").Value = 0

$ws1.Range("I7").Value = "0 de 5"
$ws1.Range("L7").Value = "0 de 5"
$ws1.Range("Q7").Value = "0 de 5"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Column widths shift (ColumnWidth uses MDW pixel-grid units that round to
# the nearest 1/MDW; add a small epsilon so the stored OOXML "width" lands
# exactly on the target integer character width instead of rounding down)
$ws2.Columns.Item(4).ColumnWidth = 9.09
$ws2.Columns.Item(5).ColumnWidth = 12.09
$ws2.Columns.Item(6).ColumnWidth = 10.09

# Month headers shift forward by one month
$ws2.Range("C1").Value = "abril"
$ws2.Range("D1").Value = "mayo"
$ws2.Range("E1").Value = "junio"
$ws2.Range("F1").Value = "julio"

# Row 2 values
$ws2.Range("C2").Value = 0
$ws2.Range("E2").Value = 2261.64
$ws2.Range("F2").Value = 0

# Row 3 values
$ws2.Range("C3").Value = 0
$ws2.Range("E3").Value = 851.4299999999999
$ws2.Range("F3").Value = 0

# Row 5 values
$ws2.Range("C5").Value = 0

# Row 6 values
$ws2.Range("C6").Value = -545.1799999999999
$ws2.Range("D6").Value = 0

# Row 7 values (totals)
$ws2.Range("C7").Value = -545.1799999999999
$ws2.Range("D7").Value = 0
$ws2.Range("E7").Value = 3113.07
$ws2.Range("F7").Value = 0
